$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the newly-added "-" placeholders in column D for the
#    "Functional difficulties" breakdown rows (child / mother sections).
$ws.Range("D28").Value = "-"
$ws.Range("D29").Value = "-"
$ws.Range("D31").Value = "-"
$ws.Range("D32").Value = "-"
$ws.Range("D33").Value = "-"

# 2. Capitalise / reword the English & Russian "by ..." section headers.
$ws.Range("C6").Value = "By sex"
$ws.Range("C12").Value = "By territory"
$ws.Range("C22").Value = "By age (in years)"
$ws.Range("C34").Value = "Wealth quintile"
$ws.Range("B22").Value = "По возрасту (в годах)"
$ws.Range("A22").Value = "Жаш курагы боюнча (жылдарда)"

# 3. Row height + wrap text for the two "functional difficulties" header
#    rows so the longer wording fits.
$ws.Range("A27:B27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 24

$ws.Range("A30:B30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 36

# 4. Update the active selection to B30.
$ws.Range("B30").Select()
